# Update NATMI ligand-receptor output with re-computed TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3845463333333334
$ws.Range("H2").Value = 1.153639
$ws.Range("I2").Value = 0.1984850200147207
$ws.Range("J2").Value = 0.1984850200147207
$ws.Range("M2").Value = 0.8229573333333334
$ws.Range("N2").Value = 2.468872
$ws.Range("O2").Value = 0.2440777672676426
$ws.Range("P2").Value = 0.2440777672676426
$ws.Range("Q2").Value = 0.3164652250231112
$ws.Range("R2").Value = 2.848187025208
$ws.Range("S2").Value = 0.04844578052126639
$ws.Range("T2").Value = 0.04844578052126638

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3845463333333334
$ws.Range("H3").Value = 1.153639
$ws.Range("I3").Value = 0.1984850200147207
$ws.Range("J3").Value = 0.1984850200147207
$ws.Range("O3").Value = 0.4345811965947162
$ws.Range("P3").Value = 0.4345811965947162
$ws.Range("Q3").Value = 0.563467282214
$ws.Range("R3").Value = 5.071205539926
$ws.Range("S3").Value = 0.08625785750412353
$ws.Range("T3").Value = 0.08625785750412353

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3845463333333334
$ws.Range("H4").Value = 1.153639
$ws.Range("I4").Value = 0.1984850200147207
$ws.Range("J4").Value = 0.1984850200147207
$ws.Range("M4").Value = 1.083466
$ws.Range("N4").Value = 3.250398
$ws.Range("O4").Value = 0.3213410361376413
$ws.Range("P4").Value = 0.3213410361376413
$ws.Range("Q4").Value = 0.4166428775913333
$ws.Range("R4").Value = 3.749785898322
$ws.Range("S4").Value = 0.06378138198933082
$ws.Range("T4").Value = 0.06378138198933082

$ws.Range("I5").Value = 0.5733580031870772
$ws.Range("J5").Value = 0.5733580031870772
$ws.Range("M5").Value = 0.8229573333333334
$ws.Range("N5").Value = 2.468872
$ws.Range("O5").Value = 0.2440777672676426
$ws.Range("P5").Value = 0.2440777672676426
$ws.Range("Q5").Value = 0.9141640486720001
$ws.Range("R5").Value = 8.227476438048001
$ws.Range("S5").Value = 0.1399439412629357
$ws.Range("T5").Value = 0.1399439412629357

$ws.Range("I6").Value = 0.5733580031870772
$ws.Range("J6").Value = 0.5733580031870772
$ws.Range("O6").Value = 0.4345811965947162
$ws.Range("P6").Value = 0.4345811965947162
$ws.Range("S6").Value = 0.2491706071021971
$ws.Range("T6").Value = 0.2491706071021971

$ws.Range("I7").Value = 0.5733580031870772
$ws.Range("J7").Value = 0.5733580031870772
$ws.Range("M7").Value = 1.083466
$ws.Range("N7").Value = 3.250398
$ws.Range("O7").Value = 0.3213410361376413
$ws.Range("P7").Value = 0.3213410361376413
$ws.Range("Q7").Value = 1.203544369848
$ws.Range("R7").Value = 10.831899328632
$ws.Range("S7").Value = 0.1842434548219444
$ws.Range("T7").Value = 0.1842434548219444

$ws.Range("G8").Value = 0.442033
$ws.Range("H8").Value = 1.326099
$ws.Range("I8").Value = 0.2281569767982021
$ws.Range("J8").Value = 0.2281569767982021
$ws.Range("M8").Value = 0.8229573333333334
$ws.Range("N8").Value = 2.468872
$ws.Range("O8").Value = 0.2440777672676426
$ws.Range("P8").Value = 0.2440777672676426
$ws.Range("Q8").Value = 0.3637742989253334
$ws.Range("R8").Value = 3.273968690328
$ws.Range("S8").Value = 0.05568804548344051
$ws.Range("T8").Value = 0.0556880454834405

$ws.Range("G9").Value = 0.442033
$ws.Range("H9").Value = 1.326099
$ws.Range("I9").Value = 0.2281569767982021
$ws.Range("J9").Value = 0.2281569767982021
$ws.Range("O9").Value = 0.4345811965947162
$ws.Range("P9").Value = 0.4345811965947162
$ws.Range("Q9").Value = 0.6477012301739998
$ws.Range("R9").Value = 5.829311071565999
$ws.Range("S9").Value = 0.09915273198839557
$ws.Range("T9").Value = 0.09915273198839557

$ws.Range("G10").Value = 0.442033
$ws.Range("H10").Value = 1.326099
$ws.Range("I10").Value = 0.2281569767982021
$ws.Range("J10").Value = 0.2281569767982021
$ws.Range("M10").Value = 1.083466
$ws.Range("N10").Value = 3.250398
$ws.Range("O10").Value = 0.3213410361376413
$ws.Range("P10").Value = 0.3213410361376413
$ws.Range("Q10").Value = 0.4789277263779999
$ws.Range("R10").Value = 4.310349537401999
$ws.Range("S10").Value = 0.07331619932636604
$ws.Range("T10").Value = 0.07331619932636606
